$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16 - LeetCode 136
$ws.Range("A16").Value = 136
$ws.Range("B16").Value = "只出现一次的数字||"
$ws.Range("D16").Value = "||"
$ws.Range("E16").Value = "位运算相关，比较新颖的做法"

# Row 17 - LeetCode 137
$ws.Range("A17").Value = 137
$ws.Range("B17").Value = "只出现一次的数字||"
$ws.Range("D17").Value = "|||"
$ws.Range("E17").Value = "位运算相关，比较新颖的做法"

# Row 18 - LeetCode 138
$ws.Range("A18").Value = 138
$ws.Range("B18").Value = "复制带随机指针的链表"
$ws.Range("D18").Value = "||"
$ws.Range("E18").Value = "背一下这个套路，套路类似于133题"

# Match the recorded selection state from the source workbook
$ws.Range("C26").Select()
